$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.531.98'
$ws.Range("E2").Value = '  +3.03%  '

$ws.Range("D3").Value = '2.627.40'
$ws.Range("E3").Value = '  +1.64%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '602.84'
$ws.Range("E5").Value = '  -0.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.24'
$ws.Range("E6").Value = '  +1.12%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").Value = '  +0.64%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.177'
$ws.Range("E9").Value = '  +5.66%  '

$ws.Range("D10").Value = '2.624.60'
$ws.Range("E10").Value = '  +1.69%  '

$ws.Range("E11").Value = '  +1.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.361'
$ws.Range("E12").Value = '  +4.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.01'
$ws.Range("E13").Value = '  -0.57%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000192'
$ws.Range("E14").Value = '  +3.70%  '

$ws.Range("D15").Value = '3.096.99'
$ws.Range("E15").Value = '  +0.66%  '

$ws.Range("D16").Value = '72.318.54'
$ws.Range("E16").Value = '  +2.98%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.64'
$ws.Range("E17").Value = '  +0.92%  '

$ws.Range("D18").Value = '2.623.36'
$ws.Range("E18").Value = '  +1.87%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.83'
$ws.Range("E19").Value = '  +5.05%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '382.38'
$ws.Range("E20").Value = '  +4.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.88'
$ws.Range("E21").Value = '  +1.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.19'
$ws.Range("E22").Value = '  +0.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.11'
$ws.Range("E23").Value = '  +12.75%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.65'
$ws.Range("E24").Value = '  +2.37%  '

$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.15%  '

$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.39'
$ws.Range("E26").Value = '  +0.77%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.13'
$ws.Range("E27").Value = '  +6.49%  '

$ws.Range("D28").Value = '2.759.73'
$ws.Range("E28").Value = '  +1.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("E30").Value = '  +2.36%  '

$ws.Range("E31").Value = '  +2.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '516.90'
$ws.Range("E32").Value = '  -1.03%  '

$ws.Range("E33").Value = '  +1.99%  '

$ws.Range("E34").Value = '  +0.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.10%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '164.65'
$ws.Range("E36").Value = '  -0.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.33'
$ws.Range("E37").Value = '  +1.42%  '

$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.113'
$ws.Range("E38").Value = '  -6.56%  '

$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.10'
$ws.Range("E39").Value = '  +0.83%  '

$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.40'
$ws.Range("E40").Value = '  +2.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.84'
$ws.Range("E41").Value = '  +2.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.09'
$ws.Range("E42").Value = '  +2.36%  '

$ws.Range("E43").Value = '  -0.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.60'
$ws.Range("E44").Value = '  +2.71%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.333'
$ws.Range("E45").Value = '  +2.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.28'
$ws.Range("E46").Value = '  +0.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '150.25'
$ws.Range("E47").Value = '  -1.64%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.70'
$ws.Range("E48").Value = '  +1.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.544'
$ws.Range("E49").Value = '  +2.89%  '

$ws.Range("E50").Value = '  +2.77%  '

$ws.Range("D51").Value = '0.0₆0262'
$ws.Range("E51").Value = '  -1.38%  '
